# Update absenteeism data rows 2-11 per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 72201
$ws.Range("B2").Value = "Felipe Carvalho"
$ws.Range("C2").Value = "Atendimento ao Cliente"
$ws.Range("D2").Value = "Consulta médica"
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 45090
$ws.Range("G2").Value = 3407.27

# Row 3
$ws.Range("A3").Value = 82666
$ws.Range("B3").Value = "Vicente Cardoso"
$ws.Range("C3").Value = "Recursos Humanos"
$ws.Range("D3").Value = "Problemas pessoais"
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 45098
$ws.Range("G3").Value = 9700.43

# Row 4
$ws.Range("A4").Value = 72237
$ws.Range("B4").Value = "Diego Gonçalves"
$ws.Range("C4").Value = "Marketing"
$ws.Range("D4").Value = "Doença"
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 45092
$ws.Range("G4").Value = 12443.69

# Row 5
$ws.Range("A5").Value = 64116
$ws.Range("B5").Value = "Giovanna Alves"
$ws.Range("C5").Value = "Vendas"
$ws.Range("D5").Value = "Problemas pessoais"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 45095
$ws.Range("G5").Value = 8716.16

# Row 6
$ws.Range("A6").Value = 79470
$ws.Range("B6").Value = "Renan da Paz"
$ws.Range("C6").Value = "Vendas"
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 45088
$ws.Range("G6").Value = 10733.85

# Row 7
$ws.Range("A7").Value = 3703
$ws.Range("B7").Value = "Gabrielly Fernandes"
$ws.Range("C7").Value = "Financeiro"
$ws.Range("D7").Value = "Outros"
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 45094
$ws.Range("G7").Value = 7258.6

# Row 8
$ws.Range("A8").Value = 92636
$ws.Range("B8").Value = "Maria Luiza Moreira"
$ws.Range("C8").Value = "TI"
$ws.Range("D8").Value = "Problemas pessoais"
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 45084
$ws.Range("G8").Value = 9300.34

# Row 9
$ws.Range("A9").Value = 50068
$ws.Range("B9").Value = "Isaac Correia"
$ws.Range("C9").Value = "Financeiro"
$ws.Range("D9").Value = "Viagem de negócios"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 45082
$ws.Range("G9").Value = 4716.16

# Row 10
$ws.Range("A10").Value = 74403
$ws.Range("B10").Value = "Dr. Luiz Fernando Rodrigues"
$ws.Range("D10").Value = "Outros"
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 45102
$ws.Range("G10").Value = 7555.19

# Row 11
$ws.Range("A11").Value = 61282
$ws.Range("B11").Value = "Sra. Elisa Campos"
$ws.Range("C11").Value = "Marketing"
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 45104
$ws.Range("G11").Value = 11573.86
